# feat: add 2022-Q3 data
#
# Target layout after edit:
#   Sheet1 "总计"     (unchanged name/position) - gets a new row for 2022-Q3,
#                      old 2020-Q4 row shifts down one row.
#   Sheet2 "2022-Q3"  (was "2020-Q4") - now holds the *new* fund-holding data,
#                      re-using the old sheet's identity (sheetId/rId).
#   Sheet3 "2020-Q4"  (new sheet)     - holds the fund-holding data that used
#                      to live in the old "2020-Q4" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Create the new sheet that will receive the *old* 2020-Q4 fund data.
#    Do this first, while $ws2 still contains that original data, so we can
#    copy it across (values + number/text types + styles) with full fidelity.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "OldQ4DataTemp"

$ws2.Range("A1:H4").Copy()
$ws3.Paste()

$ws2.Range("B1:H1").Copy()
$ws3.Range("B1:H1").PasteSpecial(-4122)

$ws2.Range("A2:A4").Copy()
$ws3.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: shift the existing 2020-Q4 row down to
#    row 3, and insert the new 2022-Q3 figures in row 2.
# ---------------------------------------------------------------------------
$ws1.Range("A2").Copy()
$ws1.Range("A3").PasteSpecial(-4122)

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2020-Q4"
$ws1.Range("C3").Value = 3
$ws1.Range("D3").Value = 0.04

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 0.03

# ---------------------------------------------------------------------------
# 3) Rename the old "2020-Q4" sheet to "2022-Q3" and overwrite its contents
#    with the new fund-holding data (reusing sheetId=2 / rId2).
# ---------------------------------------------------------------------------
$ws2.Name = "2022-Q3"

$ws2.Range("A1:H4").ClearContents()

$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"

# numeric columns A (index) and H (rank) stay plain numbers
$ws2.Range("A2").Value = 0
$ws2.Range("A3").Value = 1
$ws2.Range("A4").Value = 2
$ws2.Range("A5").Value = 3
$ws2.Range("H2").Value = 5
$ws2.Range("H3").Value = 5
$ws2.Range("H4").Value = 1
$ws2.Range("H5").Value = 1

# text columns B..G must stay text even though several look numeric
# (fund codes with leading zeros, decimal figures with trailing zeros),
# so force them to text before assigning.
$ws2.Range("B2:G5").NumberFormat = "@"

$ws2.Range("B2").Value = "014014"
$ws2.Range("C2").Value = "招商臻选平衡混合A"
$ws2.Range("D2").Value = "0.33"
$ws2.Range("E2").Value = "66.43"
$ws2.Range("F2").Value = "3.28"
$ws2.Range("G2").Value = "0.0108"

$ws2.Range("B3").Value = "014015"
$ws2.Range("C3").Value = "招商臻选平衡混合C"
$ws2.Range("D3").Value = "0.26"
$ws2.Range("E3").Value = "66.43"
$ws2.Range("F3").Value = "3.28"
$ws2.Range("G3").Value = "0.0085"

$ws2.Range("B4").Value = "011438"
$ws2.Range("C4").Value = "红塔红土盛昌优选混合A"
$ws2.Range("D4").Value = "0.14"
$ws2.Range("E4").Value = "93.85"
$ws2.Range("F4").Value = "5.58"
$ws2.Range("G4").Value = "0.0078"

$ws2.Range("B5").Value = "011439"
$ws2.Range("C5").Value = "红塔红土盛昌优选混合C"
$ws2.Range("D5").Value = "0.10"
$ws2.Range("E5").Value = "93.85"
$ws2.Range("F5").Value = "5.58"
$ws2.Range("G5").Value = "0.0056"

# Reset the "@" text formatting we applied above back to the sheet's normal
# (unstyled) look, re-using the existing default style so we don't leave a
# stray custom number format behind on these cells.
$ws1.Range("A1").Copy()
$ws2.Range("B2:G5").PasteSpecial(-4122)

# Headers (row 1) and the index column (A) use the same bold/centered style
# that is already used on the "总计" sheet (style index 2).
$ws1.Range("B1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2:A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Finish naming the sheet that holds the old 2020-Q4 fund data.
# ---------------------------------------------------------------------------
$ws3.Name = "2020-Q4"
